$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: name, email, repo link
$ws.Range("A2").Value = "عمر محمد خطاب"
$ws.Range("B2").Value = "omar564@gmail.com"
$ws.Range("C2").Value = "https://github.com/omaradds1/Open-Source-Project"

# Add a mailto hyperlink on the email cell (applies the built-in Hyperlink style)
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:omar564@gmail.com")

# Match the saved selection/active cell from the target workbook
[void]$ws.Range("C2").Select()
